$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.417.19'
$ws.Range('E2').Value = '  -0.81%  '
$ws.Range('D3').Value = '3.336.11'
$ws.Range('E3').Value = '  -4.27%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '574.54'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.39%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '178.11'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.91%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.619'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +3.40%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').Value = '3.333.79'
$ws.Range('E9').Value = '  -4.23%  '
$ws.Range('E10').Value = '  -1.99%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.85'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.08%  '
$ws.Range('E12').Value = '  -0.87%  '
$ws.Range('D13').Value = '3.912.41'
$ws.Range('E13').Value = '  -4.25%  '
$ws.Range('E14').Value = '  +0.16%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '28.41'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -4.91%  '
$ws.Range('D16').Value = '65.431.73'
$ws.Range('E16').Value = '  -0.88%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000168'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.75%  '
$ws.Range('D18').Value = '3.333.88'
$ws.Range('E18').Value = '  -4.42%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.75'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.94%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.38'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '362.56'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.14%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.43'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.86%  '
$ws.Range('E23').Value = '  +0.03%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '71.14'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.43%  '
$ws.Range('B25').Value = 'PEPE'
$ws.Range('C25').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000122'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.43%  '
$ws.Range('B26').Value = 'Polygon'
$ws.Range('C26').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.517'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.09%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.51'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.44%  '
$ws.Range('E28').Value = '  -0.91%  '
$ws.Range('E29').Value = '  -0.05%  '
$ws.Range('E30').Value = '  -1.66%  '
$ws.Range('B31').Value = 'NEARProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.60'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.87%  '
$ws.Range('B32').Value = 'USDe'
$ws.Range('C32').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.999'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.05%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '22.88'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -5.11%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.81'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.64%  '
$ws.Range('E35').Value = '  -6.18%  '
$ws.Range('E36').Value = '  -3.36%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '159.97'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.60%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.846'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -4.91%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '27.40'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -7.97%  '
$ws.Range('E40').Value = '  -0.77%  '
$ws.Range('D41').Value = '2.705.33'
$ws.Range('E41').Value = '  -4.17%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.49'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.58%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.22'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -4.39%  '
$ws.Range('E44').Value = '  -4.35%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '39.87'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.52%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0665'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.67%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '332.65'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.41%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '23.89'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.01%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0277'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.79%  '
$ws.Range('E50').Value = '  +1.86%  '
